$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dates = @(44774, 44805, 44835, 44866, 44896)
$startRow = 177

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = 269482000
}

$ws.Range("A177:A181").NumberFormat = "mmm-yy"

$ws.Range("B171").Select()
